$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.084.26"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.087.22"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "249.95"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  -5.78%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "50.93"
$ws.Range("E8").Value = "  +10.08%  "
$ws.Range("D9").Value = "60.49"
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "0.0738"
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("D13").Value = "15.14"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "2.357.85"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "0.827"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "2.091.03"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "5.05"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").Value = "37.015.18"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "72.11"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "238.76"
$ws.Range("E22").Value = "  -5.69%  "
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "2.42"
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").Value = "169.28"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").Value = "9.11"
$ws.Range("E27").Value = "  +3.14%  "
$ws.Range("D28").Value = "20.65"
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").Value = "1.99"
$ws.Range("E29").Value = "  -6.52%  "
$ws.Range("E30").Value = "  -5.34%  "
$ws.Range("D31").Value = "1.06"
$ws.Range("E31").Value = "  +20.10%  "
$ws.Range("D32").Value = "21.99"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "0.0604"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "0.0913"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +6.27%  "
$ws.Range("D38").Value = "4.07"
$ws.Range("E38").Value = "  -5.62%  "
$ws.Range("D39").Value = "1.81"
$ws.Range("E40").Value = "  -9.09%  "
$ws.Range("D41").Value = "17.72"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  +3.06%  "
$ws.Range("D44").Value = "97.74"
$ws.Range("E44").Value = "  -4.50%  "
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("D46").Value = "0.0878"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("D48").Value = "1.308.59"
$ws.Range("E48").Value = "  -3.99%  "
$ws.Range("E49").Value = "  +5.42%  "
$ws.Range("D50").Value = "2.279.14"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").Value = "  -4.14%  "
